$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.326
$ws.Range("B3").Value = 0.169
$ws.Range("B4").Value = 0.233
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0.291
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 0.293
